$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.376304149627686
$ws.Range("B1").Value = 6.287168979644775
$ws.Range("C1").Value = 7.302213668823242
$ws.Range("D1").Value = 7.415004253387451
$ws.Range("E1").Value = 3.761239051818848
